$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 54
$ws.Range("H54").Value = 5000
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 5000
$ws.Range("N54").Value = -5972

# ALC row 88
$ws.Range("H88").Value = 33336624
$ws.Range("J88").Value = 41669544
$ws.Range("L88").Value = 41669544
$ws.Range("N88").Value = -41670356

# ALC row 91
$ws.Range("H91").Value = 33336624
$ws.Range("J91").Value = 41669544
$ws.Range("L91").Value = 41669544
$ws.Range("N91").Value = -41672352

# ALC row 112
$ws.Range("H112").Value = 1381.2084
$ws.Range("J112").Value = 1492.95
$ws.Range("L112").Value = 4478.85
$ws.Range("N112").Value = -6694.85

# ALC row 129
$ws.Range("H129").Value = 545986.8
$ws.Range("I129").Value = 554.73334
$ws.Range("J129").Value = 716434.3
$ws.Range("K129").Value = 1664.20002
$ws.Range("L129").Value = 2149302.9
$ws.Range("M129").Value = 3335.79998
$ws.Range("N129").Value = -2159302.9

# ALC row 137
$ws.Range("H137").Value = 1273.3778
$ws.Range("I137").Value = 919.3143
$ws.Range("J137").Value = 2512.6
$ws.Range("K137").Value = 2757.9429
$ws.Range("L137").Value = 7537.799999999999
$ws.Range("M137").Value = -207.9429
$ws.Range("N137").Value = -12637.8

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 2446.69
$ws.Range("I32").Value = 2273.0852
$ws.Range("J32").Value = 5166.5
$ws.Range("K32").Value = 2273.0852
$ws.Range("L32").Value = 5166.5
$ws.Range("M32").Value = -1986.0852
$ws.Range("N32").Value = -5740.5

# ARM row 44
$ws.Range("H44").Value = 23200
$ws.Range("J44").Value = 23200
$ws.Range("L44").Value = 23200
$ws.Range("N44").Value = -24176

# ARM row 74
$ws.Range("H74").Value = 710
$ws.Range("I74").Value = 710
$ws.Range("K74").Value = 710
$ws.Range("M74").Value = 164

# ARM row 77
$ws.Range("H77").Value = 710
$ws.Range("I77").Value = 710
$ws.Range("K77").Value = 3550
$ws.Range("M77").Value = 818

$ws = $wb.Worksheets.Item("BSM")
# BSM row 134
$ws.Range("H134").Value = 21893.5
$ws.Range("I134").Value = 30147.742
$ws.Range("J134").Value = 2633.6
$ws.Range("K134").Value = 90443.226
$ws.Range("L134").Value = 7900.799999999999
$ws.Range("M134").Value = -87908.226
$ws.Range("N134").Value = -12970.8

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 2875533.8
$ws.Range("I31").Value = 1312.5714
$ws.Range("J31").Value = 18524072
$ws.Range("K31").Value = 1312.5714
$ws.Range("L31").Value = 18524072
$ws.Range("M31").Value = -1017.5714
$ws.Range("N31").Value = -18524662

# CRP row 34
$ws.Range("H34").Value = 2875533.8
$ws.Range("I34").Value = 1312.5714
$ws.Range("J34").Value = 18524072
$ws.Range("K34").Value = 1312.5714
$ws.Range("L34").Value = 18524072
$ws.Range("M34").Value = -1110.5714
$ws.Range("N34").Value = -18524476

# CRP row 58
$ws.Range("H58").Value = 9009925
$ws.Range("I58").Value = 938.5484
$ws.Range("K58").Value = 938.5484
$ws.Range("M58").Value = -735.5484

# CRP row 88
$ws.Range("H88").Value = 25000
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("N88").Value = -25812

# CRP row 91
$ws.Range("H91").Value = 25000
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("N91").Value = -27808

# CRP row 132
$ws.Range("H132").Value = 1894.3818
$ws.Range("I132").Value = 1767.3572
$ws.Range("K132").Value = 5302.071599999999
$ws.Range("M132").Value = -2772.071599999999

# CRP row 134
$ws.Range("H134").Value = 1314.8485
$ws.Range("I134").Value = 1263.9286
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 3791.7858
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -1256.7858
$ws.Range("N134").Value = -9870

# CRP row 136
$ws.Range("H136").Value = 9009925
$ws.Range("I136").Value = 938.5484
$ws.Range("K136").Value = 2815.6452
$ws.Range("M136").Value = -265.6451999999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 129
$ws.Range("H129").Value = 1401.9445
$ws.Range("I129").Value = 1127.2858
$ws.Range("J129").Value = 1576.7273
$ws.Range("K129").Value = 3381.8574
$ws.Range("L129").Value = 4730.1819
$ws.Range("M129").Value = 1618.1426
$ws.Range("N129").Value = -14730.1819

# CUL row 130
$ws.Range("H130").Value = 1446
$ws.Range("I130").Value = 1117.7778
$ws.Range("J130").Value = 4400
$ws.Range("K130").Value = 3353.3334
$ws.Range("L130").Value = 13200
$ws.Range("M130").Value = 1666.6666
$ws.Range("N130").Value = -23240

# CUL row 131
$ws.Range("H131").Value = 3331344.5
$ws.Range("J131").Value = 4989641.5
$ws.Range("L131").Value = 14968924.5
$ws.Range("N131").Value = -14979004.5

# CUL row 133
$ws.Range("H133").Value = 2511.0527
$ws.Range("I133").Value = 2515
$ws.Range("K133").Value = 7545
$ws.Range("M133").Value = -2485

# CUL row 134
$ws.Range("H134").Value = 2315.7144
$ws.Range("I134").Value = 2213.3333
$ws.Range("K134").Value = 6639.999899999999
$ws.Range("M134").Value = -1569.999899999999

# CUL row 136
$ws.Range("H136").Value = 1617.3572
$ws.Range("I136").Value = 1237.2727
$ws.Range("J136").Value = 3011
$ws.Range("K136").Value = 3711.8181
$ws.Range("L136").Value = 9033
$ws.Range("M136").Value = 1388.1819
$ws.Range("N136").Value = -19233

# CUL row 137
$ws.Range("H137").Value = 53475210
$ws.Range("I137").Value = 33345668
$ws.Range("J137").Value = 66056172
$ws.Range("K137").Value = 100037004
$ws.Range("L137").Value = 198168516
$ws.Range("M137").Value = -100031904
$ws.Range("N137").Value = -198178716

# CUL row 138
$ws.Range("H138").Value = 2926.6667
$ws.Range("I138").Value = 3212
$ws.Range("J138").Value = 1500
$ws.Range("K138").Value = 9636
$ws.Range("L138").Value = 4500
$ws.Range("M138").Value = -4496
$ws.Range("N138").Value = -14780

# CUL row 139
$ws.Range("H139").Value = 1964.4445
$ws.Range("I139").Value = 1964.4445
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5893.333500000001
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -753.3335000000006
$ws.Range("N139").ClearContents()

# CUL row 140
$ws.Range("H140").Value = 1675
$ws.Range("I140").Value = 1675
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 5025
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 155
$ws.Range("N140").ClearContents()

# CUL row 141
$ws.Range("H141").Value = 2461.9285
$ws.Range("I141").Value = 2372.25
$ws.Range("K141").Value = 7116.75
$ws.Range("M141").Value = -1936.75

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132
$ws.Range("H132").Value = 8335.394
$ws.Range("I132").Value = 13138.667
$ws.Range("K132").Value = 39416.001
$ws.Range("M132").Value = -36886.001

# LTW row 136
$ws.Range("H136").Value = 3172.347
$ws.Range("I136").Value = 3336.0256
$ws.Range("J136").Value = 2534
$ws.Range("K136").Value = 10008.0768
$ws.Range("L136").Value = 7602
$ws.Range("M136").Value = -7458.076799999999
$ws.Range("N136").Value = -12702

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 29161.695
$ws.Range("I122").Value = 38323.035
$ws.Range("J122").Value = 1677.6666
$ws.Range("K122").Value = 114969.105
$ws.Range("L122").Value = 5032.9998
$ws.Range("M122").Value = -112519.105
$ws.Range("N122").Value = -9932.9998

# WVR row 132
$ws.Range("H132").Value = 1456.5897
$ws.Range("I132").Value = 1265.1082
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3795.3246
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1265.3246
$ws.Range("N132").Value = -20057

# WVR row 136
$ws.Range("H136").Value = 2423.6711
$ws.Range("I136").Value = 2734.7017
$ws.Range("J136").Value = 1315.625
$ws.Range("K136").Value = 8204.105100000001
$ws.Range("L136").Value = 3946.875
$ws.Range("M136").Value = -5654.105100000001
$ws.Range("N136").Value = -9046.875
